# Update countries & provincias Spain
# Refresh the daily COVID country stats table on sheet "Pais":
#  - bump the "Datos actualizados..." timestamp
#  - update Casos totales / Nuevos casos / Casos activos / Recuperados /
#    Casos criticos / Muertes hoy / Muertes for the rows whose figures moved
#  - a few rows swap which country they show (Senegal/Tayikistan/Uzbekistan,
#    Libia/Benin, Groenlandia/Islas Malvinas, Santa Sede/Islas Turcas y
#    Caicos, Montserrat/Seychelles, Islas Virgenes Britanicas/Papua Nueva
#    Guinea) because the source list got re-sorted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 13 de Junio de 2020 a las 13:26"

# Row 7: 'India' -> 'India'
$ws.Range("B7").Value = 310131
$ws.Range("C7").Value = 528
$ws.Range("D7").Value = 154696
$ws.Range("E7").Value = 146540
$ws.Range("G7").Value = 5
$ws.Range("H7").Value = 8895

# Row 12: 'Alemania' -> 'Alemania'
$ws.Range("D12").Value = 171900
$ws.Range("E12").Value = 6488

# Row 13: 'Iran' -> 'Iran'
$ws.Range("B13").Value = 184955
$ws.Range("C13").Value = 2410
$ws.Range("D13").Value = 146748
$ws.Range("E13").Value = 29477
$ws.Range("G13").Value = 71
$ws.Range("H13").Value = 8730

# Row 23: 'Catar' -> 'Catar'
$ws.Range("B23").Value = 78416
$ws.Range("C23").Value = 1828
$ws.Range("D23").Value = 55252
$ws.Range("E23").Value = 23094

# Row 26: 'Bielorrusia' -> 'Bielorrusia'
$ws.Range("B26").Value = 53241
$ws.Range("C26").Value = 721
$ws.Range("D26").Value = 29111
$ws.Range("E26").Value = 23827
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = 303

# Row 37: 'Suiza' -> 'Suiza'
$ws.Range("B37").Value = 31094
$ws.Range("C37").Value = 31
$ws.Range("E37").Value = 356

# Row 71: 'Australia' -> 'Australia'
$ws.Range("B71").Value = 7302
$ws.Range("C71").Value = 12
$ws.Range("D71").Value = 6812
$ws.Range("E71").Value = 388

# Row 74: 'Nepal' -> 'Nepal'
$ws.Range("B74").Value = 5335
$ws.Range("C74").Value = 273
$ws.Range("D74").Value = 913
$ws.Range("E74").Value = 4404
$ws.Range("G74").Value = 2
$ws.Range("H74").Value = 18

# Row 75: 'Tayikistan' -> 'Senegal'
$ws.Range("A75").Value = "Senegal"
$ws.Range("B75").Value = 4996
$ws.Range("C75").Value = 145
$ws.Range("D75").Value = 3228
$ws.Range("E75").Value = 1708
$ws.Range("G75").Value = 4
$ws.Range("H75").Value = 60

# Row 76: 'Uzbekistan' -> 'Tayikistan'
$ws.Range("A76").Value = "Tayikistan"
$ws.Range("B76").Value = 4902
$ws.Range("C76").Value = 0
$ws.Range("D76").Value = 3158
$ws.Range("E76").Value = 1695
$ws.Range("H76").Value = 49

# Row 77: 'Senegal' -> 'Uzbekistan'
$ws.Range("A77").Value = "Uzbekistan"
$ws.Range("B77").Value = 4901
$ws.Range("C77").Value = 32
$ws.Range("D77").Value = 3758
$ws.Range("E77").Value = 1124
$ws.Range("H77").Value = 19

# Row 120: 'Madagascar' -> 'Madagascar'
$ws.Range("B120").Value = 1252
$ws.Range("C120").Value = 22
$ws.Range("D120").Value = 362
$ws.Range("E120").Value = 880

# Row 140: 'Malta' -> 'Malta'
$ws.Range("B140").Value = 646
$ws.Range("C140").Value = 1
$ws.Range("D140").Value = 601

# Row 153: 'Libia' -> 'Benin'
$ws.Range("A153").Value = "Benin"
$ws.Range("B153").Value = 412
$ws.Range("C153").Value = 24
$ws.Range("D153").Value = 222
$ws.Range("E153").Value = 184
$ws.Range("G153").Value = 1

# Row 154: 'Benin' -> 'Libia'
$ws.Range("A154").Value = "Libia"
$ws.Range("B154").Value = 409
$ws.Range("D154").Value = 59
$ws.Range("E154").Value = 344
$ws.Range("H154").Value = 6

# Row 158: 'Vietnam' -> 'Vietnam'
$ws.Range("B158").Value = 334
$ws.Range("C158").Value = 1
$ws.Range("E158").Value = 11

# Row 166: 'Gibraltar' -> 'Gibraltar'
$ws.Range("D166").Value = 173
$ws.Range("E166").Value = 3

# Row 206: 'Groenlandia' -> 'Islas Malvinas'
$ws.Range("A206").Value = "Islas Malvinas"

# Row 207: 'Islas Malvinas' -> 'Groenlandia'
$ws.Range("A207").Value = "Groenlandia"

# Row 208: 'Santa Sede' -> 'Islas Turcas y Caicos'
$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("D208").Value = 11
$ws.Range("H208").Value = 1

# Row 209: 'Islas Turcas y Caicos' -> 'Santa Sede'
$ws.Range("A209").Value = "Santa Sede"
$ws.Range("D209").Value = 12
$ws.Range("H209").Value = 0

# Row 210: 'Montserrat' -> 'Seychelles'
$ws.Range("A210").Value = "Seychelles"
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

# Row 211: 'Seychelles' -> 'Montserrat'
$ws.Range("A211").Value = "Montserrat"
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1

# Row 213: 'Islas Virgenes Britanicas' -> 'Papua Nueva Guinea'
$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("D213").Value = 8
$ws.Range("H213").Value = 0

# Row 214: 'Papua Nueva Guinea' -> 'Islas Virgenes Britanicas'
$ws.Range("A214").Value = "Islas Virgenes Britanicas"
$ws.Range("D214").Value = 7
$ws.Range("H214").Value = 1
